# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 203
    3  = 1029
    6  = 4583
    7  = 26
    8  = 380
    9  = 1338
    10 = 886
    11 = 54
    12 = 956
    14 = 535
    16 = 250
    17 = 22
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
